# Applies the "Added content to Readme" commit:
#  1. Slide 13 - update the GitHub handle shown in the contact info.
#  2. Slide 3  - collapse the "url + handle + ';'" run triples for the
#                Rocketseat / LinkedIn / GitHub social links into single
#                runs, updating the GitHub handle and dropping the
#                trailing ';' on the GitHub line.
#  3. Slide 4  - append two new paragraphs pointing at the project's
#                GitHub repository, the second one being a bulleted
#                hyperlink.

$p = $ppt.ActivePresentation

function Replace-Substring {
    param($TextRange, [string]$Needle, [string]$NewText)

    $full = $TextRange.Text
    $idx = $full.IndexOf($Needle)
    if ($idx -lt 0) {
        throw "Needle not found: $Needle"
    }
    $sub = $TextRange.Characters($idx + 1, $Needle.Length)
    $sub.Text = $NewText
}

# ---------------------------------------------------------------------
# 1. Slide 13 ("Obrigado!") - github.com/PedroLacombe -> .../pedrohenriquelacombe
# ---------------------------------------------------------------------
$slide13 = $p.Slides.Item(13)
$shape13 = $slide13.Shapes.Item(2)
$tr13 = $shape13.TextFrame.TextRange
Replace-Substring $tr13 "https://github.com/PedroLacombe" "https://github.com/pedrohenriquelacombe"

# ---------------------------------------------------------------------
# 2. Slide 3 ("Sobre mim...") - merge the split runs for the social links
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$shape3 = $slide3.Shapes.Item(2)
$tr3 = $shape3.TextFrame.TextRange

Replace-Substring $tr3 "https://app.rocketseat.com.br/me/fuskinha;" "https://app.rocketseat.com.br/me/fuskinha"
Replace-Substring $tr3 "https://www.linkedin.com/in/pedrohenriquelacombe;" "https://www.linkedin.com/in/pedrohenriquelacombe"
Replace-Substring $tr3 "https://github.com/PedroLacombe;" "https://github.com/pedrohenriquelacombe"

# ---------------------------------------------------------------------
# 3. Slide 4 ("Sobre a aula...") - add the GitHub repository paragraphs
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$shape4 = $slide4.Shapes.Item(2)
$tr4 = $shape4.TextFrame.TextRange

# First new paragraph: plain heading line, no bullet.
$tr4.InsertAfter([char]13 + "Repositório do projeto no GitHub:") | Out-Null

$full4 = $tr4.Text
$headIdx = $full4.IndexOf("Repositório do projeto no GitHub:")
$headRange = $tr4.Characters($headIdx + 1, "Repositório do projeto no GitHub:".Length)
$headRange.ParagraphFormat.Bullet.Visible = $false

# Second new paragraph: bulleted sub-level link to the repository.
$repoUrl = "https://github.com/rocketseat-experts-club/spring-cloud-openfeign-2021-07-03"
$tr4.InsertAfter([char]13 + $repoUrl) | Out-Null

$full4b = $tr4.Text
$urlIdx = $full4b.IndexOf($repoUrl)
$urlRange = $tr4.Characters($urlIdx + 1, $repoUrl.Length)
$urlRange.IndentLevel = 2
$urlRange.ActionSettings.Item(1).Hyperlink.Address = $repoUrl
